# Fix bug in ApplicationPropertiesMiddleware test data: three expected-value
# cells in the ScopedTraceLogger TestJson sheet were missing the "header*"
# prefix that the middleware actually emits. Correct them to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G14").Value = "header*X-User=curly&X-Set-ScopedLogger=curly"
$ws.Range("G16").Value = "header*X-User=shemp"
$ws.Range("G18").Value = "header*X-User=curly&X-Clear-ScopedLogger=curly"

# Restore the view the author left the sheet in (scrolled down/right,
# cursor resting on the last corrected cell).
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("G19").Select() | Out-Null
